$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Data")

# Update the existing week-ending 11/1 (row 34) figures: final revised Positive count
$ws.Range("B34").Value = 327

# Add the new week-ending 11/7 row (row 35)
$ws.Range("A35").Value = 44142
$ws.Range("B35").Value = 286
$ws.Range("C35").Formula = "=AVERAGE(B32:B35)"
$ws.Range("D35").Formula = "=(B35/126884)*100000"
$ws.Range("E35").Value = 60
$ws.Range("F35").Formula = "=F34+B35"
$ws.Range("G35").Formula = "=G34+E35"

# Match the author's updated selection state (active cell moves to the new last row)
$null = $ws.Range("F35").Select()
